$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard_Tests")

# ------------------------------------------------------------------
# 1. Widen column C to fit the new, longer test-step text.
# ------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 107.5546875

# ------------------------------------------------------------------
# 2. Make room: the sheet grows from 17 to 23 rows. Rows 16 and 17
#    (Verify Logout Functionality / Verify URL contains "login")
#    slide down to rows 19 and 20, three brand-new test cases are
#    inserted above them (rows 16-18) and three blank rows are
#    appended at the bottom (21-23).
#
#    Blanket-format rows 14-23 first (border box, no fill - the
#    same formatting every data row in this sheet already uses),
#    then patch the handful of cells that need a different style.
# ------------------------------------------------------------------
foreach ($r in 14..23) {
    $ws.Range("A16:J16").Copy() | Out-Null
    $ws.Range("A${r}:J${r}").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# Column C on the "description" rows uses the shaded/boxed style (3).
foreach ($r in 14,15,16,18,19,20) {
    $ws.Range("C10").Copy() | Out-Null
    $ws.Range("C${r}").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# B14 keeps the box-without-bottom-border style (11) used elsewhere.
$ws.Range("A8").Copy() | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# C17 loses its border formatting entirely (plain/general style).
$ws.Range("A6").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3. Push the existing "Logout" test case down to rows 19-20 (its new
#    test-case id is applied afterwards, once the new rows above it
#    have claimed their shared-string slots).
# ------------------------------------------------------------------
$ws.Range("B19").Value = $ws.Range("B16").Value2
$ws.Range("C19").Value = $ws.Range("C16").Value2
$ws.Range("C20").Value = $ws.Range("C17").Value2

# ------------------------------------------------------------------
# 4. Fill in the three new "live dashboard" database test cases.
#    Values are entered in this specific order so new shared strings
#    line up the same way the source workbook produced them.
# ------------------------------------------------------------------
$ws.Range("B16").Value = "Verify Live Customer Count"
$ws.Range("A18").Value = "TC_DB_08"
$ws.Range("C16").Value = "1.Verify text ""{DB_QUERY}SELECT COUNT(*) FROM customers"" at ""//h5[text()='Total Customers']/following-sibling::h2"""
$ws.Range("A17").Value = "TC_DB_07"
$ws.Range("B17").Value = "Verify Live Product Count"
$ws.Range("C17").Value = "1.Verify text ""{DB_QUERY}SELECT COUNT(*) FROM products"" at ""//h5[text()='Total Products']/following-sibling::h2"""
$ws.Range("B18").Value = "Verify Live Revenue"
$ws.Range("C18").Value = "1.Verify text ""{DB_QUERY}SELECT SUM(total_bill) FROM sales"" at ""//h5[text()='Total Revenue']/following-sibling::h2"""
$ws.Range("A19").Value = "TC_DB_09"

# ------------------------------------------------------------------
# 5. Final selection, matching where the author ended up.
# ------------------------------------------------------------------
$ws.Range("C23").Select()
